$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCellValue($Worksheet, $CellRef, $Text, $Row) {
    $target = $Worksheet.Range($CellRef)
    $donor = $Worksheet.Cells.Item($Row, 2)
    $target.NumberFormat = "@"
    $target.Value = $Text
    $target.Style = $donor.Style
}

Set-TextCellValue $ws "D2" "25.776.97" 2
Set-TextCellValue $ws "E2" "  -0.16%  " 2
Set-TextCellValue $ws "D3" "1.635.45" 3
Set-TextCellValue $ws "E3" "  +0.06%  " 3
Set-TextCellValue $ws "E4" "  -0.16%  " 4
Set-TextCellValue $ws "D5" "215.74" 5
Set-TextCellValue $ws "E5" "  +0.36%  " 5
Set-TextCellValue $ws "E7" "  -0.16%  " 7
Set-TextCellValue $ws "E8" "  -0.06%  " 8
Set-TextCellValue $ws "E9" "  -1.13%  " 9
Set-TextCellValue $ws "D10" "19.54" 10
Set-TextCellValue $ws "E10" "  -1.94%  " 10
Set-TextCellValue $ws "D11" "0.0791" 11
Set-TextCellValue $ws "E11" "  +1.54%  " 11
Set-TextCellValue $ws "D12" "4.25" 12
Set-TextCellValue $ws "E12" "  +0.14%  " 12
Set-TextCellValue $ws "D13" "1.861.25" 13
Set-TextCellValue $ws "E13" "  +0.07%  " 13
Set-TextCellValue $ws "D14" "1.638.83" 14
Set-TextCellValue $ws "E14" "  +0.22%  " 14
Set-TextCellValue $ws "D15" "0.563" 15
Set-TextCellValue $ws "E15" "  +0.68%  " 15
Set-TextCellValue $ws "D16" "0.0₃0764" 16
Set-TextCellValue $ws "E16" "  -0.56%  " 16
Set-TextCellValue $ws "D17" "63.23" 17
Set-TextCellValue $ws "E17" "  +0.30%  " 17
Set-TextCellValue $ws "D18" "25.808.30" 18
Set-TextCellValue $ws "E18" "  -0.08%  " 18
Set-TextCellValue $ws "D21" "192.41" 21
Set-TextCellValue $ws "E21" "  -0.74%  " 21
Set-TextCellValue $ws "E22" "  +0.58%  " 22
Set-TextCellValue $ws "D23" "6.31" 23
Set-TextCellValue $ws "E23" "  +2.42%  " 23
Set-TextCellValue $ws "D24" "1.85" 24
Set-TextCellValue $ws "E24" "  +4.99%  " 24
Set-TextCellValue $ws "E25" "  -0.10%  " 25
Set-TextCellValue $ws "D26" "141.79" 26
Set-TextCellValue $ws "E26" "  +1.75%  " 26
Set-TextCellValue $ws "E27" "  +1.62%  " 27
Set-TextCellValue $ws "E28" "  +1.03%  " 28
Set-TextCellValue $ws "D29" "15.50" 29
Set-TextCellValue $ws "E29" "  +0.29%  " 29
Set-TextCellValue $ws "E30" "  -0.13%  " 30
Set-TextCellValue $ws "E31" "  -0.37%  " 31
Set-TextCellValue $ws "E33" "  -0.63%  " 33
Set-TextCellValue $ws "E34" "  -0.54%  " 34
Set-TextCellValue $ws "E35" "  -0.34%  " 35
Set-TextCellValue $ws "E36" "  +0.28%  " 36
Set-TextCellValue $ws "D37" "1.131.16" 37
Set-TextCellValue $ws "E37" "  +1.36%  " 37
Set-TextCellValue $ws "E38" "  -2.03%  " 38
Set-TextCellValue $ws "E39" "  -0.91%  " 39
Set-TextCellValue $ws "D40" "0.0155" 40
Set-TextCellValue $ws "E40" "  -0.84%  " 40
Set-TextCellValue $ws "E41" "  +0.09%  " 41
Set-TextCellValue $ws "D42" "2.53" 42
Set-TextCellValue $ws "E42" "  +0.73%  " 42
Set-TextCellValue $ws "E43" "  +0.31%  " 43
Set-TextCellValue $ws "D44" "100.71" 44
Set-TextCellValue $ws "E44" "  +1.40%  " 44
Set-TextCellValue $ws "D45" "0.802" 45
Set-TextCellValue $ws "E45" "  +0.33%  " 45
Set-TextCellValue $ws "D46" "1.770.44" 46
Set-TextCellValue $ws "E46" "  -0.13%  " 46
Set-TextCellValue $ws "E47" "  +2.18%  " 47
Set-TextCellValue $ws "E48" "  -0.27%  " 48
Set-TextCellValue $ws "E49" "  -0.95%  " 49
Set-TextCellValue $ws "E50" "  -0.21%  " 50
Set-TextCellValue $ws "E51" "  +4.22%  " 51
